# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto symbol
# list with the latest scraped values. The sheet stores these as plain text
# (e.g. "279.45", "6.89%"), so each cell's NumberFormat is forced to "@"
# (Text) before the assignment to stop Excel from auto-converting the
# numeric-looking/percent-looking strings into real numbers, then the style
# is reset back to "Normal" so no residual text-format styling is left
# behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (values scraped from crypto price feed refresh)
$updates = @{
    'D2' = '279.45'
    'E2' = '6.89%'
    'D3' = '27.28'
    'E3' = '2.02%'
    'D4' = '4.824'
    'E4' = '2.62%'
    'D5' = '0.06333'
    'E5' = '2.19%'
    'D6' = '6.936'
    'E6' = '2.79%'
    'D7' = '3.384'
    'E7' = '6.69%'
    'D8' = '0.8787'
    'E8' = '3.54%'
    'D9' = '0.9493'
    'E9' = '4.18%'
    'D10' = '0.1469'
    'E10' = '4.60%'
    'D11' = '0.05143'
    'E11' = '1.30%'
    'D12' = '0.07276'
    'E12' = '2.34%'
    'D13' = '0.03147'
    'E13' = '1.98%'
    'D14' = '0.09074'
    'E14' = '0.34%'
    'D15' = '0.001552'
    'E15' = '0.57%'
    'D16' = '0.0006281'
    'E16' = '1.61%'
    'D17' = '0.005913'
    'E17' = '-0.56%'
    'D18' = '3.444'
    'E18' = '-0.15%'
    'D19' = '2.292'
    'E19' = '4.77%'
    'D20' = '0.3161'
    'E20' = '1.86%'
    'D21' = '0.1311'
    'E21' = '0.07%'
    'D22' = '3.878'
    'E22' = '-5.55%'
    'D23' = '0.04324'
    'E23' = '1.76%'
    'D24' = '0.001176'
    'E24' = '-0.04%'
    'D25' = '0.004294'
    'E25' = '5.79%'
    'D26' = '0.0001189'
    'E26' = '-0.98%'
    'D27' = '0.0001689'
    'E27' = '2.99%'
    'D40' = '0.04083'
    'E40' = '3.09%'
    'D41' = '0.006679'
    'E41' = '61.55%'
    'D42' = '0.1162'
    'E42' = '4.51%'
    'D43' = '0.002199'
    'E43' = '3.06%'
    'D44' = '0.01316'
    'E44' = '-0.81%'
    'D45' = '0.00005217'
    'E45' = '1.03%'
    'D46' = '0.00000000749'
    'E46' = '-0.16%'
    'E47' = '852.58%'
    'D48' = '0.02249'
    'E48' = '-33.91%'
    'D49' = '0.00002099'
    'E49' = '-0.16%'
    'D50' = '0.0001999'
    'E50' = '-0.16%'
}

foreach ($cellRef in $updates.Keys) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cellRef]
    $rng.Style = "Normal"
}

